$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update headline metrics after trade #190 closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(3, 2).Value = 1399.95   # Current Capital
$summary.Cells.Item(4, 2).Value = -0.27     # Total P&L $
$summary.Cells.Item(5, 2).Value = -0.02     # Total P&L %
$summary.Cells.Item(6, 2).Value = 218       # Total Trades
$summary.Cells.Item(7, 2).Value = 84        # Winning Trades
$summary.Cells.Item(9, 2).Value = 38.53     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(5, 3).Value = 99.95   # Capital
$status.Cells.Item(5, 4).Value = 185     # Trades
$status.Cells.Item(5, 5).Value = -0.38   # P&L $
$status.Cells.Item(5, 6).Value = -0.05   # P&L %
$status.Cells.Item(5, 7).Value = 37.84   # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "All Trades": close out trade #218 (row 219) + append rows for
# trades #251 (volatility_scorer) and #252 (MarketMaking)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(219, 7).Value = 0.141408                 # Exit Price
$allTrades.Cells.Item(219, 8).Value = "CLOSED"                  # Status
$allTrades.Cells.Item(219, 9).Value = 8.775700000000001         # P&L %
$allTrades.Cells.Item(219, 10).Value = 0.01                     # P&L $
$allTrades.Cells.Item(219, 11).Value = 99.95                    # Capital After
$allTrades.Cells.Item(219, 12).Value = "early_exit"              # Exit Reason
$allTrades.Cells.Item(219, 13).Value = 0.12                     # Duration (min)

# New row 252 -> trade #251 (volatility_scorer, still OPEN)
$allTrades.Cells.Item(252, 1).Value = 251
$allTrades.Cells.Item(252, 2).NumberFormat = "@"
$allTrades.Cells.Item(252, 2).Value = "2026-02-17"
$allTrades.Cells.Item(252, 2).Style = "Normal"
$allTrades.Cells.Item(252, 3).Value = "22:07:33"
$allTrades.Cells.Item(252, 4).Value = "volatility_scorer"
$allTrades.Cells.Item(252, 5).Value = "NEUTRAL"
$allTrades.Cells.Item(252, 6).Value = 0.13
$allTrades.Cells.Item(252, 8).Value = "OPEN"
$allTrades.Cells.Item(252, 9).Value = 0
$allTrades.Cells.Item(252, 10).Value = 0
$allTrades.Cells.Item(252, 11).Value = 100
$allTrades.Cells.Item(252, 13).Value = 0
$allTrades.Cells.Item(252, 14).Value = 0
$allTrades.Cells.Item(252, 15).Value = 0
$allTrades.Cells.Item(252, 16).Value = 0.85
$allTrades.Cells.Item(252, 17).Value = "Low vol market (score: inf) - ideal for market making"

# New row 253 -> trade #252 (MarketMaking, still OPEN)
$allTrades.Cells.Item(253, 1).Value = 252
$allTrades.Cells.Item(253, 2).NumberFormat = "@"
$allTrades.Cells.Item(253, 2).Value = "2026-02-17"
$allTrades.Cells.Item(253, 2).Style = "Normal"
$allTrades.Cells.Item(253, 3).Value = "22:07:33"
$allTrades.Cells.Item(253, 4).Value = "MarketMaking"
$allTrades.Cells.Item(253, 5).Value = "UP"
$allTrades.Cells.Item(253, 6).Value = 0.83
$allTrades.Cells.Item(253, 8).Value = "OPEN"
$allTrades.Cells.Item(253, 9).Value = 0
$allTrades.Cells.Item(253, 10).Value = 0
$allTrades.Cells.Item(253, 11).Value = 99.93696837561131
$allTrades.Cells.Item(253, 13).Value = 0
$allTrades.Cells.Item(253, 14).Value = 0
$allTrades.Cells.Item(253, 15).Value = 0
$allTrades.Cells.Item(253, 16).Value = 0.6
$allTrades.Cells.Item(253, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# Sheet "volatility_scorer": append row for trade #251
# ---------------------------------------------------------------------------
$volScorer = $wb.Worksheets.Item("volatility_scorer")
$volScorer.Cells.Item(8, 1).Value = 251
$volScorer.Cells.Item(8, 2).NumberFormat = "@"
$volScorer.Cells.Item(8, 2).Value = "2026-02-17"
$volScorer.Cells.Item(8, 2).Style = "Normal"
$volScorer.Cells.Item(8, 3).Value = "22:07:33"
$volScorer.Cells.Item(8, 4).Value = "volatility_scorer"
$volScorer.Cells.Item(8, 5).Value = "NEUTRAL"
$volScorer.Cells.Item(8, 6).Value = 0.13
$volScorer.Cells.Item(8, 8).Value = "OPEN"
$volScorer.Cells.Item(8, 9).Value = 0
$volScorer.Cells.Item(8, 10).Value = 0
$volScorer.Cells.Item(8, 11).Value = 100
$volScorer.Cells.Item(8, 12).Value = 0
$volScorer.Cells.Item(8, 13).Value = 0
$volScorer.Cells.Item(8, 14).Value = 0.85
$volScorer.Cells.Item(8, 15).Value = "Low vol market (score: inf) - ideal for market making"
$volScorer.Cells.Item(8, 17).Value = 0

# ---------------------------------------------------------------------------
# Sheet "MarketMaking": close out trade #218 (row 186) + append row for
# trade #252
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$marketMaking.Cells.Item(186, 7).Value = 0.141408                # Exit Price
$marketMaking.Cells.Item(186, 8).Value = "CLOSED"                 # Status
$marketMaking.Cells.Item(186, 9).Value = 8.775700000000001        # P&L %
$marketMaking.Cells.Item(186, 10).Value = 0.01                    # P&L $
$marketMaking.Cells.Item(186, 11).Value = 99.95                   # Capital After
$marketMaking.Cells.Item(186, 16).Value = "early_exit"             # Exit Reason
$marketMaking.Cells.Item(186, 17).Value = 0.12                    # Duration (min)

# New row 213 -> trade #252
$marketMaking.Cells.Item(213, 1).Value = 252
$marketMaking.Cells.Item(213, 2).NumberFormat = "@"
$marketMaking.Cells.Item(213, 2).Value = "2026-02-17"
$marketMaking.Cells.Item(213, 2).Style = "Normal"
$marketMaking.Cells.Item(213, 3).Value = "22:07:33"
$marketMaking.Cells.Item(213, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(213, 5).Value = "UP"
$marketMaking.Cells.Item(213, 6).Value = 0.83
$marketMaking.Cells.Item(213, 8).Value = "OPEN"
$marketMaking.Cells.Item(213, 9).Value = 0
$marketMaking.Cells.Item(213, 10).Value = 0
$marketMaking.Cells.Item(213, 11).Value = 99.93696837561131
$marketMaking.Cells.Item(213, 12).Value = 0
$marketMaking.Cells.Item(213, 13).Value = 0
$marketMaking.Cells.Item(213, 14).Value = 0.6
$marketMaking.Cells.Item(213, 15).Value = "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item(213, 17).Value = 0
